# Actualización automática 2025-06-20 14:55:09
#
# Updates ALMEIDA CUATIN JHONATHANN CARLOS / COMFALASDI COMPAÑIA FAMILIAR
# LASCANO DIAZ C. LTDA. junio sales figures for "240X80 PORCELANATO" and
# "PORCELANATO", plus all the dependent summary / totals / percentages
# that are derived from them across the three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-product-group sales for the client row.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D8").Value = 982.46      # 240X80 PORCELANATO
$wsGrupo.Range("M8").Value = 782.58      # PORCELANATO

# M column non-zero counter row ("n de 28") - M8 became non-zero so the
# count goes up by one, from 4 to 5.
$wsGrupo.Range("M30").Value = "5 de 28"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": monthly (junio) totals for the client row, and
# the column total in row 30.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F8").Value = 1765.04
$wsMensual.Range("F30").Value = 5567.02

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO per
# product group, plus the TOTAL row.
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 982.46
$wsCumpl.Range("E3").Value = 2137.6545
$wsCumpl.Range("F3").Value = 0.3148794699681695

# Row 16: PORCELANATO
$wsCumpl.Range("D16").Value = 3539.89
$wsCumpl.Range("E16").Value = 15258.72
$wsCumpl.Range("F16").Value = 0.1883059438969158

# Row 19: TOTAL
$wsCumpl.Range("D19").Value = 5561.26
$wsCumpl.Range("E19").Value = 23976.53107555787
$wsCumpl.Range("F19").Value = 0.1882760964005148

# Column F width narrowed slightly (25 -> 24 characters).
$wsCumpl.Columns.Item(6).ColumnWidth = 23.14
